$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ROW50-FE-LIFTER": append new row 97 (dimension A1:I96 -> A1:I97)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 97
$ws.Cells.Item($r, 1).Value = 45772.8092425
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x3a"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 314
$ws.Cells.Item($r, 9).Value = 14

# ---------------------------------------------------------------------------
# Sheet "ROW50-MID-LIFTER": append new row 99 (dimension A1:I98 -> A1:I99)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 99
$ws.Cells.Item($r, 1).Value = 45772.77182870371
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x3e"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$gCell = $ws.Cells.Item($r, 7)
$gCell.NumberFormat = "@"
$gCell.Value = "568631262647113771663628"
$gCell.Style = "Normal"
$ws.Cells.Item($r, 8).Value = 318
$ws.Cells.Item($r, 9).Value = 25

# ---------------------------------------------------------------------------
# Sheet "ROW11-FE-LIFTER": append new row 97 (dimension A1:I96 -> A1:I97)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 97
$ws.Cells.Item($r, 1).Value = 45772.84066738426
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x3a"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 314
$ws.Cells.Item($r, 9).Value = 20

# ---------------------------------------------------------------------------
# Sheet "ROW11-MID-LIFTER": append new row 97 (dimension A1:I96 -> A1:I97)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 97
$ws.Cells.Item($r, 1).Value = 45772.95915216435
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x42"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 322
$ws.Cells.Item($r, 9).Value = 25
